$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(4).Insert()

Write-Output "done"
